$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-08-31 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-01 Friday", 2) | Out-Null
$d.Content.Find.Execute("27×82=2214", $true, $false, $false, $false, $false, $true, 1, $false, "88×81=7128", 2) | Out-Null
$d.Content.Find.Execute("74×36=2664", $true, $false, $false, $false, $false, $true, 1, $false, "94×80=7520", 2) | Out-Null
$d.Content.Find.Execute("69×61=4209", $true, $false, $false, $false, $false, $true, 1, $false, "43×25=1075", 2) | Out-Null
$d.Content.Find.Execute("25×70=1750", $true, $false, $false, $false, $false, $true, 1, $false, "80×14=1120", 2) | Out-Null
$d.Content.Find.Execute("26×60=1560", $true, $false, $false, $false, $false, $true, 1, $false, "59×51=3009", 2) | Out-Null
$d.Content.Find.Execute("71×60=4260", $true, $false, $false, $false, $false, $true, 1, $false, "39×60=2340", 2) | Out-Null
$d.Content.Find.Execute("33×36=1188", $true, $false, $false, $false, $false, $true, 1, $false, "61×31=1891", 2) | Out-Null
$d.Content.Find.Execute("97×44=4268", $true, $false, $false, $false, $false, $true, 1, $false, "53×55=2915", 2) | Out-Null
$d.Content.Find.Execute("31×23=713", $true, $false, $false, $false, $false, $true, 1, $false, "76×57=4332", 2) | Out-Null
$d.Content.Find.Execute("61×55=3355", $true, $false, $false, $false, $false, $true, 1, $false, "66×58=3828", 2) | Out-Null
$d.Content.Find.Execute("44×18=792", $true, $false, $false, $false, $false, $true, 1, $false, "18×35=630", 2) | Out-Null
$d.Content.Find.Execute("68×63=4284", $true, $false, $false, $false, $false, $true, 1, $false, "67×77=5159", 2) | Out-Null
$d.Content.Find.Execute("80×46=3680", $true, $false, $false, $false, $false, $true, 1, $false, "67×36=2412", 2) | Out-Null
$d.Content.Find.Execute("98×39=3822", $true, $false, $false, $false, $false, $true, 1, $false, "53×64=3392", 2) | Out-Null
$d.Content.Find.Execute("36×30=1080", $true, $false, $false, $false, $false, $true, 1, $false, "14×15=210", 2) | Out-Null
$d.Content.Find.Execute("99×75=7425", $true, $false, $false, $false, $false, $true, 1, $false, "54×17=918", 2) | Out-Null
$d.Content.Find.Execute("12×37=444", $true, $false, $false, $false, $false, $true, 1, $false, "42×43=1806", 2) | Out-Null
$d.Content.Find.Execute("72×49=3528", $true, $false, $false, $false, $false, $true, 1, $false, "91×87=7917", 2) | Out-Null
$d.Content.Find.Execute("43×78=3354", $true, $false, $false, $false, $false, $true, 1, $false, "44×90=3960", 2) | Out-Null
$d.Content.Find.Execute("81×85=6885", $true, $false, $false, $false, $false, $true, 1, $false, "69×81=5589", 2) | Out-Null
$d.Content.Find.Execute("59×45=2655", $true, $false, $false, $false, $false, $true, 1, $false, "12×34=408", 2) | Out-Null
$d.Content.Find.Execute("99×63=6237", $true, $false, $false, $false, $false, $true, 1, $false, "24×61=1464", 2) | Out-Null
$d.Content.Find.Execute("45×97=4365", $true, $false, $false, $false, $false, $true, 1, $false, "80×39=3120", 2) | Out-Null
$d.Content.Find.Execute("58×27=1566", $true, $false, $false, $false, $false, $true, 1, $false, "39×93=3627", 2) | Out-Null
$d.Content.Find.Execute("36×77=2772", $true, $false, $false, $false, $false, $true, 1, $false, "61×60=3660", 2) | Out-Null
